$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the placeholder "waffles" nutrition-label values with "placeholder"
# for each menu item row.
$ws.Range("F2").Value = "placeholder"
$ws.Range("F3").Value = "placeholder"
$ws.Range("F4").Value = "placeholder"

# Shrink the table (and its autofilter) so it no longer includes the
# trailing empty row that had no data.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:G4"))

# Update the active selection to the empty row right below the table.
$ws.Rows.Item(5).Select() | Out-Null
